# doc: update upload file
# Rewrites the "Parts", "Category" and "Webshop" sheets with the refreshed
# computer-parts catalogue, clears stale rows, fixes up selections/active
# sheet, and strips the now-unused cell style from the Webshop sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Category sheet: was A1:A5, becomes A1:A12 with the expanded category
# list. Filled in first (matching the original authoring order: the new
# category rows were added before the Parts sheet referenced them).
# ---------------------------------------------------------------------
$wsCategory = $wb.Worksheets.Item("Category")

$wsCategory.Range("A1").Value = "CategoryName"
$wsCategory.Range("A2").Value = "Monitor"
$wsCategory.Range("A6").Value = "Motherboard"
$wsCategory.Range("A7").Value = "Memory "
$wsCategory.Range("A8").Value = "Power Supplie"
$wsCategory.Range("A9").Value = "Computer Case"
$wsCategory.Range("A10").Value = "Graphics Card"
$wsCategory.Range("A11").Value = "Network Card"
$wsCategory.Range("A4").Value = "Processors"
$wsCategory.Range("A12").Value = "Sound Card"
$wsCategory.Range("A3").Value = "hard drive "
$wsCategory.Range("A5").Value = "DVD"

# Selection on Category moves to a full-row selection of row 6.
$wsCategory.Rows.Item(6).Select()

# ---------------------------------------------------------------------
# Parts sheet: was A1:B11, becomes A1:B8 with new header/content pairs.
# ---------------------------------------------------------------------
$wsParts = $wb.Worksheets.Item("Parts")

# Clear the three rows that no longer exist so the used range shrinks to
# A1:B8 (ClearContents removes the row entries entirely).
$wsParts.Range("A9:B11").ClearContents()

$wsParts.Range("A1").Value = "ComputerPartName"
$wsParts.Range("B1").Value = "CategoryName"

$wsParts.Range("A2").Value = "24 Inch Full HD Monitor"
$wsParts.Range("B2").Value = "Monitor"

$wsParts.Range("A3").Value = "1TB HDD"
$wsParts.Range("B3").Value = "hard drive "

$wsParts.Range("A4").Value = "Intel Core i5"
$wsParts.Range("B4").Value = "Processors"

$wsParts.Range("A5").Value = "motherboard atx DDR4"
$wsParts.Range("B5").Value = "Motherboard"

$wsParts.Range("A6").Value = "ddr4 64GB"
$wsParts.Range("B6").Value = "Memory "

$wsParts.Range("A7").Value = "ATX Mid Tower Gaming Case"
$wsParts.Range("B7").Value = "Computer Case"

$wsParts.Range("A8").Value = "AMD Radeon RX"
$wsParts.Range("B8").Value = "Graphics Card"

# ---------------------------------------------------------------------
# Webshop sheet: values are re-ordered/re-indexed (same logical content)
# and the two cells that previously had an explicit "applyFill" style
# lose that formatting.
# ---------------------------------------------------------------------
$wsWebshop = $wb.Worksheets.Item("Webshop")

$wsWebshop.Range("B5:B6").ClearFormats()

$wsWebshop.Range("A1").Value = "WebshopName"
$wsWebshop.Range("B1").Value = "WebshopURL"

$wsWebshop.Range("A2").Value = "Amazon"
$wsWebshop.Range("B2").Value = "Amazon.com"

$wsWebshop.Range("A3").Value = "EBAY"
$wsWebshop.Range("B3").Value = "Ebay.com"

$wsWebshop.Range("A4").Value = "Currys"
$wsWebshop.Range("B4").Value = "https://www.currys.co.uk/"

$wsWebshop.Range("A5").Value = "hobbycraft"
$wsWebshop.Range("B5").Value = "https://www.hobbycraft.co.uk/"

$wsWebshop.Range("A6").Value = "hm"
$wsWebshop.Range("B6").Value = "hm.com"

# ---------------------------------------------------------------------
# Active sheet / selections: the workbook now opens on "Parts" (tab 0),
# which was previously "Webshop" (tab 2). Selecting on Webshop first
# (it keeps the A7 selection it already had, just without being the
# active tab any more), then activating "Parts" and selecting A7 there
# last, so "Parts" ends up as the active/tabSelected sheet.
# ---------------------------------------------------------------------
$wsWebshop.Range("A7").Select()

$wsParts.Activate()
$wsParts.Range("A7").Select()
